$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row=2; A='wenden'; B='none'; C='none' },
    @{ Row=3; A='posten'; B='face/face023.jpg'; C='face' },
    @{ Row=4; A='landen'; B='flower/flower032.jpg'; C='flower' },
    @{ Row=5; A='rufen'; B='none'; C='none' },
    @{ Row=6; A='lehnen'; B='face/face006.jpg'; C='face' },
    @{ Row=7; A='decken'; B='face/face001.jpg'; C='face' },
    @{ Row=8; A='enden'; B='none'; C='none' },
    @{ Row=9; A='triefen'; B='face/face028.jpg'; C='face' },
    @{ Row=10; A='deuten'; B='face/face019.jpg'; C='face' },
    @{ Row=11; A='opfern'; B='none'; C='none' },
    @{ Row=12; A='bieten'; B='face/face031.jpg'; C='face' },
    @{ Row=13; A='tropfen'; B='flower/flower025.jpg'; C='flower' },
    @{ Row=14; A='schalten'; B='none'; C='none' },
    @{ Row=15; A='ärgern'; B='flower/flower013.jpg'; C='flower' },
    @{ Row=16; A='seufzen'; B='face/face015.jpg'; C='face' },
    @{ Row=17; A='klagen'; B='none'; C='none' },
    @{ Row=18; A='dringen'; B='face/face014.jpg'; C='face' },
    @{ Row=19; A='kehren'; B='flower/flower017.jpg'; C='flower' },
    @{ Row=20; A='sparen'; B='none'; C='none' },
    @{ Row=21; A='stürmen'; B='flower/flower014.jpg'; C='flower' },
    @{ Row=22; A='nerven'; B='face/face007.jpg'; C='face' },
    @{ Row=23; A='hören'; B='none'; C='none' },
    @{ Row=24; A='runden'; B='face/face021.jpg'; C='face' },
    @{ Row=25; A='herrschen'; B='flower/flower012.jpg'; C='flower' },
    @{ Row=26; A='drohen'; B='none'; C='none' },
    @{ Row=27; A='starren'; B='flower/flower028.jpg'; C='flower' },
    @{ Row=28; A='schreiben'; B='flower/flower023.jpg'; C='flower' },
    @{ Row=29; A='orten'; B='none'; C='none' },
    @{ Row=30; A='steuern'; B='face/face020.jpg'; C='face' },
    @{ Row=31; A='beißen'; B='flower/flower015.jpg'; C='flower' },
    @{ Row=32; A='weigern'; B='none'; C='none' },
    @{ Row=33; A='ehren'; B='flower/flower002.jpg'; C='flower' },
    @{ Row=34; A='retten'; B='face/face018.jpg'; C='face' },
    @{ Row=35; A='ächzen'; B='none'; C='none' },
    @{ Row=36; A='reizen'; B='face/face012.jpg'; C='face' },
    @{ Row=37; A='binden'; B='face/face029.jpg'; C='face' },
    @{ Row=38; A='kosten'; B='none'; C='none' },
    @{ Row=39; A='bauen'; B='flower/flower020.jpg'; C='flower' },
    @{ Row=40; A='schleppen'; B='flower/flower005.jpg'; C='flower' },
    @{ Row=41; A='dauern'; B='none'; C='none' },
    @{ Row=42; A='schwimmen'; B='face/face017.jpg'; C='face' },
    @{ Row=43; A='streichen'; B='flower/flower003.jpg'; C='flower' },
    @{ Row=44; A='stören'; B='none'; C='none' },
    @{ Row=45; A='albern'; B='flower/flower024.jpg'; C='flower' },
    @{ Row=46; A='achten'; B='face/face013.jpg'; C='face' },
    @{ Row=47; A='bremsen'; B='none'; C='none' },
    @{ Row=48; A='zögern'; B='flower/flower019.jpg'; C='flower' },
    @{ Row=49; A='parken'; B='flower/flower029.jpg'; C='flower' }
)

# Clear the existing data rows so stale shared-string entries are dropped
$ws.Range("A2:C49").ClearContents()

# Re-populate column by column (A, then B, then C) to rebuild the shared-string table
# in the same traversal order used by the source edit.
foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.A
}
foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 2).Value = $r.B
}
foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 3).Value = $r.C
}
